$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 340, pushing the existing rows 340:361 down to 341:362.
$ws.Rows.Item(340).Insert()

# Populate the new row 340 with a new weekly record, using the same
# static/categorical values as the surrounding rows (copied from row 341,
# which now holds the data that used to be in row 340).
$ws.Range("A340").Value2 = $ws.Range("A341").Value2
$ws.Range("B340").Value2 = $ws.Range("B341").Value2
$ws.Range("C340").Value2 = $ws.Range("C341").Value2
$ws.Range("D340").Value2 = 44746
$ws.Range("E340").Value2 = $ws.Range("E341").Value2
$ws.Range("F340").Value2 = $ws.Range("F341").Value2
$ws.Range("G340").Value2 = $ws.Range("G341").Value2
$ws.Range("H340").Value2 = $ws.Range("H341").Value2
$ws.Range("I340").Value2 = $ws.Range("I341").Value2
$ws.Range("J340").Value2 = 172
$ws.Range("K340").Value2 = 4000
$ws.Range("L340").Value2 = 4500
$ws.Range("M340").Value2 = 4253
$ws.Range("N340").Value2 = $ws.Range("N341").Value2
$ws.Range("O340").Value2 = $ws.Range("O341").Value2
$ws.Range("P340").Value2 = 1418
$ws.Range("Q340").Value2 = $ws.Range("Q341").Value2
$ws.Range("R340").Value2 = $ws.Range("R341").Value2
